$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("clause")

# Copy the formatting of the last existing data row (49) down onto the
# three new rows (50-52) so the new rows pick up the same cell styles
# (font/number format) used by the rest of the table.
$ws.Range("A49:I49").Copy()
$ws.Range("A50:I52").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 50 - beilu / sz123082
$ws.Range("A50").Value = "beilu"
$ws.Range("B50").Value = "sz123082"
$ws.Range("C50").Value = 0.5
$ws.Range("D50").Value = 0.7
$ws.Range("E50").Value = 1.2
$ws.Range("F50").Value = 1.8
$ws.Range("G50").Value = 2.5
$ws.Range("H50").Value = 15
$ws.Range("I50").Value = 46362

# Row 51 - ruida / sz128116
$ws.Range("A51").Value = "ruida"
$ws.Range("B51").Value = "sz128116"
$ws.Range("C51").Value = 0.4
$ws.Range("D51").Value = 0.5
$ws.Range("E51").Value = 1
$ws.Range("F51").Value = 1.5
$ws.Range("G51").Value = 1.9
$ws.Range("H51").Value = 10
$ws.Range("I51").Value = 46201

# Row 52 - huazheng / sh113639
$ws.Range("A52").Value = "huazheng"
$ws.Range("B52").Value = "sh113639"
$ws.Range("C52").Value = 0.2
$ws.Range("D52").Value = 0.4
$ws.Range("E52").Value = 0.6
$ws.Range("F52").Value = 1.5
$ws.Range("G52").Value = 1.8
$ws.Range("H52").Value = 8
$ws.Range("I52").Value = 46775

# Widen column B a bit (matches the new, longer codes now shown there).
$ws.Columns.Item(2).ColumnWidth = 10.75

# Scroll the view down and select the last-entered cells, as a user would
# after adding the new rows.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D52:E52").Select()
